# Auto-generated edit script: update computed profit columns (H-N) on 36 rows
# across all 8 sheets, per the scheduled-runner price refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(18, 8).Value = 200200660
$ws.Cells.Item(18, 9).Value = 333666750
$ws.Cells.Item(18, 10).Value = 1475
$ws.Cells.Item(18, 11).Value = 333666750
$ws.Cells.Item(18, 12).Value = 1475
$ws.Cells.Item(18, 13).Value = -333666466
$ws.Cells.Item(18, 14).Value = -2043
$ws.Cells.Item(28, 8).Value = 1146.5
$ws.Cells.Item(28, 9).Value = 1146.5
$ws.Cells.Item(28, 10).Value = 0
$ws.Cells.Item(28, 11).Value = 1146.5
$ws.Cells.Item(28, 12).Value = 0
$ws.Cells.Item(28, 13).Value = -661.5
$ws.Cells.Item(28, 14).ClearContents()
$ws.Cells.Item(53, 8).Value = 184
$ws.Cells.Item(53, 9).Value = 179.5
$ws.Cells.Item(53, 10).Value = 187.6
$ws.Cells.Item(53, 11).Value = 179.5
$ws.Cells.Item(53, 12).Value = 187.6
$ws.Cells.Item(53, 13).Value = 457.5
$ws.Cells.Item(53, 14).Value = -1461.6
$ws.Cells.Item(74, 8).Value = 3413.3333
$ws.Cells.Item(74, 9).Value = 0
$ws.Cells.Item(74, 10).Value = 3413.3333
$ws.Cells.Item(74, 11).Value = 0
$ws.Cells.Item(74, 12).Value = 3413.3333
$ws.Cells.Item(74, 13).ClearContents()
$ws.Cells.Item(74, 14).Value = -5285.3333
$ws.Cells.Item(77, 8).Value = 3413.3333
$ws.Cells.Item(77, 9).Value = 0
$ws.Cells.Item(77, 10).Value = 3413.3333
$ws.Cells.Item(77, 11).Value = 0
$ws.Cells.Item(77, 12).Value = 17066.6665
$ws.Cells.Item(77, 13).ClearContents()
$ws.Cells.Item(77, 14).Value = -26426.6665
$ws.Cells.Item(86, 8).Value = 1870.3
$ws.Cells.Item(86, 9).Value = 1800.4
$ws.Cells.Item(86, 10).Value = 2080
$ws.Cells.Item(86, 11).Value = 1800.4
$ws.Cells.Item(86, 12).Value = 2080
$ws.Cells.Item(86, 13).Value = -677.4000000000001
$ws.Cells.Item(86, 14).Value = -4326
$ws.Cells.Item(89, 8).Value = 1870.3
$ws.Cells.Item(89, 9).Value = 1800.4
$ws.Cells.Item(89, 10).Value = 2080
$ws.Cells.Item(89, 11).Value = 9002
$ws.Cells.Item(89, 12).Value = 10400
$ws.Cells.Item(89, 13).Value = -3386
$ws.Cells.Item(89, 14).Value = -21632
$ws.Cells.Item(107, 8).Value = 1197.4445
$ws.Cells.Item(107, 9).Value = 1103.8667
$ws.Cells.Item(107, 10).Value = 1665.3334
$ws.Cells.Item(107, 11).Value = 1103.8667
$ws.Cells.Item(107, 12).Value = 1665.3334
$ws.Cells.Item(107, 13).Value = 816.1333
$ws.Cells.Item(107, 14).Value = -5505.3334
$ws.Cells.Item(113, 8).Value = 1998.75
$ws.Cells.Item(113, 10).Value = 1998.75
$ws.Cells.Item(113, 12).Value = 1998.75
$ws.Cells.Item(113, 14).Value = -8506.75
$ws.Cells.Item(135, 8).Value = 1145
$ws.Cells.Item(135, 9).Value = 823.9
$ws.Cells.Item(135, 10).Value = 1787.2
$ws.Cells.Item(135, 11).Value = 7415.099999999999
$ws.Cells.Item(135, 12).Value = 16084.8
$ws.Cells.Item(135, 13).Value = -4880.099999999999
$ws.Cells.Item(135, 14).Value = -21154.8
$ws.Cells.Item(141, 8).Value = 1200
$ws.Cells.Item(141, 9).Value = 1097.5
$ws.Cells.Item(141, 10).Value = 1302.5
$ws.Cells.Item(141, 11).Value = 3292.5
$ws.Cells.Item(141, 12).Value = 3907.5
$ws.Cells.Item(141, 13).Value = 1887.5
$ws.Cells.Item(141, 14).Value = -14267.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 1598
$ws.Cells.Item(2, 9).Value = 1653.1538
$ws.Cells.Item(2, 11).Value = 1653.1538
$ws.Cells.Item(2, 13).Value = -1540.1538
$ws.Cells.Item(15, 8).Value = 5000
$ws.Cells.Item(15, 10).Value = 5000
$ws.Cells.Item(15, 12).Value = 5000
$ws.Cells.Item(15, 14).Value = -5700
$ws.Cells.Item(45, 8).Value = 3535
$ws.Cells.Item(45, 9).Value = 2956
$ws.Cells.Item(45, 10).Value = 4114
$ws.Cells.Item(45, 11).Value = 2956
$ws.Cells.Item(45, 12).Value = 4114
$ws.Cells.Item(45, 13).Value = -2579
$ws.Cells.Item(45, 14).Value = -4868
$ws.Cells.Item(56, 8).Value = 12000
$ws.Cells.Item(56, 10).Value = 12000
$ws.Cells.Item(56, 12).Value = 12000
$ws.Cells.Item(56, 14).Value = -13484
$ws.Cells.Item(74, 8).Value = 3493.1794
$ws.Cells.Item(74, 9).Value = 4159.3335
$ws.Cells.Item(74, 10).Value = 1272.6666
$ws.Cells.Item(74, 11).Value = 4159.3335
$ws.Cells.Item(74, 12).Value = 1272.6666
$ws.Cells.Item(74, 13).Value = -3285.3335
$ws.Cells.Item(74, 14).Value = -3020.6666
$ws.Cells.Item(77, 8).Value = 3493.1794
$ws.Cells.Item(77, 9).Value = 4159.3335
$ws.Cells.Item(77, 10).Value = 1272.6666
$ws.Cells.Item(77, 11).Value = 20796.6675
$ws.Cells.Item(77, 12).Value = 6363.333000000001
$ws.Cells.Item(77, 13).Value = -16428.6675
$ws.Cells.Item(77, 14).Value = -15099.333
$ws.Cells.Item(116, 8).Value = 1598
$ws.Cells.Item(116, 9).Value = 1653.1538
$ws.Cells.Item(116, 11).Value = 1653.1538
$ws.Cells.Item(116, 13).Value = 640.8462

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 1598
$ws.Cells.Item(3, 9).Value = 1653.1538
$ws.Cells.Item(3, 11).Value = 1653.1538
$ws.Cells.Item(3, 13).Value = -1539.1538

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(107, 8).Value = 953.6667
$ws.Cells.Item(107, 9).Value = 988.0909
$ws.Cells.Item(107, 10).Value = 924.53845
$ws.Cells.Item(107, 11).Value = 988.0909
$ws.Cells.Item(107, 12).Value = 924.53845
$ws.Cells.Item(107, 13).Value = 931.9091
$ws.Cells.Item(107, 14).Value = -4764.53845

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(12, 8).Value = 45.75
$ws.Cells.Item(12, 9).Value = 38
$ws.Cells.Item(12, 10).Value = 48.333332
$ws.Cells.Item(12, 11).Value = 114
$ws.Cells.Item(12, 12).Value = 144.999996
$ws.Cells.Item(12, 13).Value = 59
$ws.Cells.Item(12, 14).Value = -490.999996

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(15, 8).Value = 14450
$ws.Cells.Item(15, 10).Value = 14450
$ws.Cells.Item(15, 12).Value = 14450
$ws.Cells.Item(15, 14).Value = -15026
$ws.Cells.Item(80, 8).Value = 2990.0688
$ws.Cells.Item(80, 10).Value = 3193.2307
$ws.Cells.Item(80, 12).Value = 3193.2307
$ws.Cells.Item(80, 14).Value = -5189.2307
$ws.Cells.Item(81, 8).Value = 14450
$ws.Cells.Item(81, 10).Value = 14450
$ws.Cells.Item(81, 12).Value = 14450
$ws.Cells.Item(81, 14).Value = -16446
$ws.Cells.Item(83, 8).Value = 2990.0688
$ws.Cells.Item(83, 10).Value = 3193.2307
$ws.Cells.Item(83, 12).Value = 15966.1535
$ws.Cells.Item(83, 14).Value = -25950.1535
$ws.Cells.Item(84, 8).Value = 14450
$ws.Cells.Item(84, 10).Value = 14450
$ws.Cells.Item(84, 12).Value = 43350
$ws.Cells.Item(84, 14).Value = -53334

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(36, 8).Value = 74333.336
$ws.Cells.Item(36, 10).Value = 74333.336
$ws.Cells.Item(36, 12).Value = 74333.336
$ws.Cells.Item(36, 14).Value = -75457.336
$ws.Cells.Item(56, 8).Value = 6316.6665
$ws.Cells.Item(56, 9).Value = 4975
$ws.Cells.Item(56, 10).Value = 9000
$ws.Cells.Item(56, 11).Value = 4975
$ws.Cells.Item(56, 12).Value = 9000
$ws.Cells.Item(56, 13).Value = -4284
$ws.Cells.Item(56, 14).Value = -10382
$ws.Cells.Item(61, 8).Value = 15964.286
$ws.Cells.Item(61, 9).Value = 15964.286
$ws.Cells.Item(61, 10).Value = 0
$ws.Cells.Item(61, 11).Value = 15964.286
$ws.Cells.Item(61, 12).Value = 0
$ws.Cells.Item(61, 13).Value = -15762.286
$ws.Cells.Item(61, 14).ClearContents()
$ws.Cells.Item(82, 8).Value = 0
$ws.Cells.Item(82, 9).Value = 0
$ws.Cells.Item(82, 10).Value = 0
$ws.Cells.Item(82, 11).Value = 0
$ws.Cells.Item(82, 12).Value = 0
$ws.Cells.Item(82, 13).ClearContents()
$ws.Cells.Item(82, 14).ClearContents()
$ws.Cells.Item(85, 8).Value = 0
$ws.Cells.Item(85, 9).Value = 0
$ws.Cells.Item(85, 10).Value = 0
$ws.Cells.Item(85, 11).Value = 0
$ws.Cells.Item(85, 12).Value = 0
$ws.Cells.Item(85, 13).ClearContents()
$ws.Cells.Item(85, 14).ClearContents()
$ws.Cells.Item(113, 8).Value = 15964.286
$ws.Cells.Item(113, 9).Value = 15964.286
$ws.Cells.Item(113, 10).Value = 0
$ws.Cells.Item(113, 11).Value = 15964.286
$ws.Cells.Item(113, 12).Value = 0
$ws.Cells.Item(113, 13).Value = -13794.286
$ws.Cells.Item(113, 14).ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(49, 8).Value = 10018.333
$ws.Cells.Item(49, 9).Value = 56
$ws.Cells.Item(49, 10).Value = 14999.5
$ws.Cells.Item(49, 11).Value = 56
$ws.Cells.Item(49, 12).Value = 14999.5
$ws.Cells.Item(49, 13).Value = 174
$ws.Cells.Item(49, 14).Value = -15459.5
$ws.Cells.Item(94, 8).Value = 15249.5
$ws.Cells.Item(94, 10).Value = 15249.5
$ws.Cells.Item(94, 12).Value = 15249.5
$ws.Cells.Item(94, 14).Value = -17051.5
$ws.Cells.Item(107, 8).Value = 361.2143
$ws.Cells.Item(107, 9).Value = 294.8
$ws.Cells.Item(107, 10).Value = 527.25
$ws.Cells.Item(107, 11).Value = 884.4000000000001
$ws.Cells.Item(107, 12).Value = 1581.75
$ws.Cells.Item(107, 13).Value = 1035.6
$ws.Cells.Item(107, 14).Value = -5421.75
$ws.Cells.Item(113, 8).Value = 334.68182
$ws.Cells.Item(113, 9).Value = 322.5238
$ws.Cells.Item(113, 10).Value = 590
$ws.Cells.Item(113, 11).Value = 967.5714
$ws.Cells.Item(113, 12).Value = 1770
$ws.Cells.Item(113, 13).Value = 1202.4286
$ws.Cells.Item(113, 14).Value = -6110
